$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
$ws.Range('D2').Value = '30.894.18'
$ws.Range('E2').Value = '  +0.63%  '
# Row 3
$ws.Range('D3').Value = '1.923.51'
# Row 4
Set-TextValue 'D4' '1.003'
$ws.Range('E4').Value = '  +0.11%  '
# Row 5
Set-TextValue 'D5' '240.06'
$ws.Range('E5').Value = '  -3.05%  '
# Row 6
$ws.Range('E6').Value = '  +0.16%  '
# Row 7
Set-TextValue 'D7' '0.4909'
$ws.Range('E7').Value = '  -0.54%  '
# Row 8
Set-TextValue 'D8' '0.2973'
$ws.Range('E8').Value = '  +0.57%  '
# Row 9
Set-TextValue 'D9' '0.06775'
$ws.Range('E9').Value = '  -0.61%  '
# Row 10
$ws.Range('D10').Value = '1.944.23'
$ws.Range('E10').Value = '  +2.87%  '
# Row 11
Set-TextValue 'D11' '17.08'
$ws.Range('E11').Value = '  -0.63%  '
# Row 12
Set-TextValue 'D12' '0.07310'
$ws.Range('E12').Value = '  +1.03%  '
# Row 13
Set-TextValue 'D13' '5.176'
$ws.Range('E13').Value = '  +2.12%  '
# Row 14
Set-TextValue 'D14' '89.59'
$ws.Range('E14').Value = '  -2.37%  '
# Row 15
Set-TextValue 'D15' '0.6723'
$ws.Range('E15').Value = '  -0.74%  '
# Row 16
$ws.Range('D16').Value = '30.862.62'
$ws.Range('E16').Value = '  +0.60%  '
# Row 17
Set-TextValue 'D17' '0.000008007'
$ws.Range('E17').Value = '  +0.46%  '
# Row 18
Set-TextValue 'D18' '13.54'
$ws.Range('E18').Value = '  +2.44%  '
# Row 19
$ws.Range('E19').Value = '  +0.17%  '
# Row 20
$ws.Range('D20').Value = '2.171.31'
$ws.Range('E20').Value = '  +1.59%  '
# Row 21
Set-TextValue 'D21' '1.003'
$ws.Range('E21').Value = '  +0.00%  '
# Row 22
Set-TextValue 'D22' '5.192'
$ws.Range('E22').Value = '  +7.47%  '
# Row 23
Set-TextValue 'D23' '206.15'
$ws.Range('E23').Value = '  +9.62%  '
# Row 24
Set-TextValue 'D24' '6.297'
$ws.Range('E24').Value = '  +3.97%  '
# Row 25
Set-TextValue 'D25' '9.665'
# Row 26
Set-TextValue 'D26' '160.62'
$ws.Range('E26').Value = '  +2.51%  '
# Row 27
Set-TextValue 'D27' '19.05'
$ws.Range('E27').Value = '  -0.10%  '
# Row 28
Set-TextValue 'D28' '1.985'
$ws.Range('E28').Value = '  +4.04%  '
# Row 29
Set-TextValue 'D29' '1.430'
$ws.Range('E29').Value = '  +1.85%  '
# Row 30
Set-TextValue 'D30' '4.369'
$ws.Range('E30').Value = '  +1.51%  '
# Row 31
Set-TextValue 'D31' '0.09204'
$ws.Range('E31').Value = '  +2.31%  '
# Row 32
Set-TextValue 'D32' '4.068'
$ws.Range('E32').Value = '  +1.45%  '
# Row 33
Set-TextValue 'D33' '0.05199'
$ws.Range('E33').Value = '  +0.33%  '
# Row 34
Set-TextValue 'D34' '0.7545'
$ws.Range('E34').Value = '  +1.29%  '
# Row 35
Set-TextValue 'D35' '1.123'
# Row 36
Set-TextValue 'D36' '2.726'
$ws.Range('E36').Value = '  -0.13%  '
# Row 37
Set-TextValue 'D37' '0.01861'
$ws.Range('E37').Value = '  +1.13%  '
# Row 38
Set-TextValue 'D38' '2.725'
$ws.Range('E38').Value = '  +2.03%  '
# Row 39
Set-TextValue 'D39' '0.9268'
$ws.Range('E39').Value = '  -1.39%  '
# Row 40
Set-TextValue 'D40' '2.089'
$ws.Range('E40').Value = '  -3.27%  '
# Row 41
Set-TextValue 'D41' '0.4513'
$ws.Range('E41').Value = '  +2.08%  '
# Row 42
Set-TextValue 'D42' '108.34'
$ws.Range('E42').Value = '  +2.74%  '
# Row 43
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D43' '72.28'
$ws.Range('E43').Value = '  +24.96%  '
# Row 44
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D44' '5.928'
$ws.Range('E44').Value = '  +2.87%  '
# Row 45
$ws.Range('E45').Value = '  +1.08%  '
# Row 46
Set-TextValue 'D46' '0.1397'
$ws.Range('E46').Value = '  +4.28%  '
# Row 47
Set-TextValue 'D47' '7.701'
$ws.Range('E47').Value = '  +0.97%  '
# Row 48
Set-TextValue 'D48' '36.05'
$ws.Range('E48').Value = '  +7.71%  '
# Row 49
Set-TextValue 'D49' '9.052'
$ws.Range('E49').Value = '  +4.20%  '
# Row 50
Set-TextValue 'D50' '0.05949'
$ws.Range('E50').Value = '  +1.80%  '
# Row 51
Set-TextValue 'D51' '0.4086'
$ws.Range('E51').Value = '  +3.68%  '
